$wb = $excel.ActiveWorkbook

# ---- Sheet: Overview ----
$wsOverview = $wb.Worksheets.Item("Overview")

# Set plain (non-hyperlink) cell values
$wsOverview.Range("A1").Value = "File Name"
$wsOverview.Range("B1").Value = "zh-cn"
$wsOverview.Range("C1").Value = "de-de"
$wsOverview.Range("D1").Value = "Latest Handoff Date"
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsOverview.Range("D2").Value = "2016-03-24 15:23:58"
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-03-24 15:23:58"
$wsOverview.Range("B4").Value = "Ready for handoff"
$wsOverview.Range("C4").Value = "Ready for handoff"
$wsOverview.Range("D4").Value = "2016-03-24 15:23:58"
$wsOverview.Range("B5").Value = "Ready for handoff"
$wsOverview.Range("C5").Value = "Ready for handoff"
$wsOverview.Range("D5").Value = "2016-03-24 15:23:58"

# Reset hyperlinks (clearing then re-adding in final order)
$wsOverview.Range("A1").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/18c195304a5ec01160e997ad2750ee74e1688391/e2e/calleeMd1.md", [Type]::Missing, [Type]::Missing, "calleeMd1.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/18c195304a5ec01160e997ad2750ee74e1688391/e2e/calleeMd2.md", [Type]::Missing, [Type]::Missing, "calleeMd2.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/18c195304a5ec01160e997ad2750ee74e1688391/e2e/callerMd1.md", [Type]::Missing, [Type]::Missing, "callerMd1.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/18c195304a5ec01160e997ad2750ee74e1688391/e2e/callerMd2.md", [Type]::Missing, [Type]::Missing, "callerMd2.md") | Out-Null

# ---- Sheet: zh-cn ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# Set plain (non-hyperlink) cell values
$wsZhCn.Range("A1").Value = "Source File Name"
$wsZhCn.Range("B1").Value = "File Extension"
$wsZhCn.Range("C1").Value = "Status"
$wsZhCn.Range("D1").Value = "Latest Handoff File"
$wsZhCn.Range("E1").Value = "Latest Handoff Datetime"
$wsZhCn.Range("F1").Value = "Latest Target File"
$wsZhCn.Range("G1").Value = "Latest Handback File"
$wsZhCn.Range("H1").Value = "Latest Handback DateTime"
$wsZhCn.Range("I1").Value = "Reference Tokens"
$wsZhCn.Range("J1").Value = "Handoff Reason"
$wsZhCn.Range("K1").Value = "Dependency From"
$wsZhCn.Range("L1").Value = "Error Detail"
$wsZhCn.Range("B2").Value = ".md"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("E2").Value = "2016-03-24 15:23:53"
$wsZhCn.Range("H2").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("J2").Value = "Include"
$wsZhCn.Range("K2").Value = "e2e\callerMd2.md,`ne2e\callerMd1.md"
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "2016-03-24 15:23:53"
$wsZhCn.Range("H3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("J3").Value = "Include"
$wsZhCn.Range("K3").Value = "e2e\callerMd1.md"
$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = "Ready for handoff"
$wsZhCn.Range("E4").Value = "2016-03-24 15:23:53"
$wsZhCn.Range("H4").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("I4").Value = "e2e\calleeMd1.md,`ne2e\calleeMd2.md"
$wsZhCn.Range("J4").Value = "Include"
$wsZhCn.Range("B5").Value = ".md"
$wsZhCn.Range("C5").Value = "Ready for handoff"
$wsZhCn.Range("E5").Value = "2016-03-24 15:23:53"
$wsZhCn.Range("H5").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("I5").Value = "e2e\calleeMd1.md"
$wsZhCn.Range("J5").Value = "Include"

# Reset hyperlinks (clearing then re-adding in final order)
$wsZhCn.Range("A1").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/18c195304a5ec01160e997ad2750ee74e1688391/e2e/calleeMd1.md", [Type]::Missing, [Type]::Missing, "calleeMd1.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/15d95caaa5fb180759aa8f3c85e7f1a396fdf471/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.zh-cn.xlf") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/18c195304a5ec01160e997ad2750ee74e1688391/e2e/calleeMd2.md", [Type]::Missing, [Type]::Missing, "calleeMd2.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/15d95caaa5fb180759aa8f3c85e7f1a396fdf471/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.zh-cn.xlf") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/18c195304a5ec01160e997ad2750ee74e1688391/e2e/callerMd1.md", [Type]::Missing, [Type]::Missing, "callerMd1.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/15d95caaa5fb180759aa8f3c85e7f1a396fdf471/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.zh-cn.xlf") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/18c195304a5ec01160e997ad2750ee74e1688391/e2e/callerMd2.md", [Type]::Missing, [Type]::Missing, "callerMd2.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/15d95caaa5fb180759aa8f3c85e7f1a396fdf471/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.zh-cn.xlf") | Out-Null

# ---- Sheet: de-de ----
$wsDeDe = $wb.Worksheets.Item("de-de")

# Set plain (non-hyperlink) cell values
$wsDeDe.Range("A1").Value = "Source File Name"
$wsDeDe.Range("B1").Value = "File Extension"
$wsDeDe.Range("C1").Value = "Status"
$wsDeDe.Range("D1").Value = "Latest Handoff File"
$wsDeDe.Range("E1").Value = "Latest Handoff Datetime"
$wsDeDe.Range("F1").Value = "Latest Target File"
$wsDeDe.Range("G1").Value = "Latest Handback File"
$wsDeDe.Range("H1").Value = "Latest Handback DateTime"
$wsDeDe.Range("I1").Value = "Reference Tokens"
$wsDeDe.Range("J1").Value = "Handoff Reason"
$wsDeDe.Range("K1").Value = "Dependency From"
$wsDeDe.Range("L1").Value = "Error Detail"
$wsDeDe.Range("B2").Value = ".md"
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("E2").Value = "2016-03-24 15:23:58"
$wsDeDe.Range("H2").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("J2").Value = "Include"
$wsDeDe.Range("K2").Value = "e2e\callerMd2.md,`ne2e\callerMd1.md"
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "2016-03-24 15:23:58"
$wsDeDe.Range("H3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("J3").Value = "Include"
$wsDeDe.Range("K3").Value = "e2e\callerMd1.md"
$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = "Ready for handoff"
$wsDeDe.Range("E4").Value = "2016-03-24 15:23:58"
$wsDeDe.Range("H4").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("I4").Value = "e2e\calleeMd1.md,`ne2e\calleeMd2.md"
$wsDeDe.Range("J4").Value = "Include"
$wsDeDe.Range("B5").Value = ".md"
$wsDeDe.Range("C5").Value = "Ready for handoff"
$wsDeDe.Range("E5").Value = "2016-03-24 15:23:58"
$wsDeDe.Range("H5").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("I5").Value = "e2e\calleeMd1.md"
$wsDeDe.Range("J5").Value = "Include"

# Reset hyperlinks (clearing then re-adding in final order)
$wsDeDe.Range("A1").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/18c195304a5ec01160e997ad2750ee74e1688391/e2e/calleeMd1.md", [Type]::Missing, [Type]::Missing, "calleeMd1.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e336536ab78b5628f3a68039e5d8e0b49d49faa3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.de-de.xlf", [Type]::Missing, [Type]::Missing, "calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.de-de.xlf") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/18c195304a5ec01160e997ad2750ee74e1688391/e2e/calleeMd2.md", [Type]::Missing, [Type]::Missing, "calleeMd2.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e336536ab78b5628f3a68039e5d8e0b49d49faa3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.de-de.xlf", [Type]::Missing, [Type]::Missing, "calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.de-de.xlf") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/18c195304a5ec01160e997ad2750ee74e1688391/e2e/callerMd1.md", [Type]::Missing, [Type]::Missing, "callerMd1.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e336536ab78b5628f3a68039e5d8e0b49d49faa3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.de-de.xlf", [Type]::Missing, [Type]::Missing, "callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.de-de.xlf") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/18c195304a5ec01160e997ad2750ee74e1688391/e2e/callerMd2.md", [Type]::Missing, [Type]::Missing, "callerMd2.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e336536ab78b5628f3a68039e5d8e0b49d49faa3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.de-de.xlf", [Type]::Missing, [Type]::Missing, "callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.de-de.xlf") | Out-Null

